$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2,5).Value = 24
$ws.Cells.Item(2,6).Value = 8
# Row 3
$ws.Cells.Item(3,5).Value = 18
$ws.Cells.Item(3,6).Value = 3
# Row 4
$ws.Cells.Item(4,5).Value = 16
$ws.Cells.Item(4,6).Value = 2
# Row 5
$ws.Cells.Item(5,6).Value = 7
# Row 6
$ws.Cells.Item(6,1).Value = "R"
$ws.Cells.Item(6,2).Value = "perceived risks"
$ws.Cells.Item(6,3).Value = "Material"
$ws.Cells.Item(6,4).Value = "Addresses physical and structural considerations, such as weaknesses in soft robots or design flaws impacting resilience and performance."
$ws.Cells.Item(6,5).Value = 9
$ws.Cells.Item(6,6).Value = 1
# Row 7
$ws.Cells.Item(7,1).Value = "TP"
$ws.Cells.Item(7,2).Value = "perceived technological possibilities"
$ws.Cells.Item(7,3).Value = "Accessibility"
$ws.Cells.Item(7,4).Value = "Identifies instances where participants highlight the advantage of robots, particularly drones, in accessing and navigating locations that are challenging for humans to reach, including unstable areas."
$ws.Cells.Item(7,6).Value = 22
# Row 8
$ws.Cells.Item(8,3).Value = "Reliability"
$ws.Cells.Item(8,4).Value = "Emphasizes consistent performance, precise control, and the role of rescue robots in reducing risks by ensuring dependable operation in hazardous conditions."
$ws.Cells.Item(8,6).Value = 23
# Row 9
$ws.Cells.Item(9,1).Value = "SA"
$ws.Cells.Item(9,2).Value = "perceived safety"
$ws.Cells.Item(9,3).Value = "Access"
$ws.Cells.Item(9,4).Value = "Focuses on the enhanced ability of rescue robots to access remote or hard-to-reach locations, enabling quicker response times and potentially saving lives."
$ws.Cells.Item(9,6).Value = 20
# Row 10
$ws.Cells.Item(10,3).Value = "Physical Capabilities"
$ws.Cells.Item(10,4).Value = "Discusses the ability of rescue robots to perform physical tasks beyond human capabilities, emphasizing their potential in challenging environments."
$ws.Cells.Item(10,5).Value = 6
$ws.Cells.Item(10,6).Value = 22
# Row 11
$ws.Cells.Item(11,3).Value = "Delivery of Goods"
$ws.Cells.Item(11,4).Value = "Highlights the potential use of rescue robots for delivering goods in hazardous environments, like war zones, emphasizing the safety benefits."
$ws.Cells.Item(11,5).Value = 6
$ws.Cells.Item(11,6).Value = 8
# Row 12
$ws.Cells.Item(12,1).Value = "SA"
$ws.Cells.Item(12,2).Value = "perceived safety"
$ws.Cells.Item(12,3).Value = "Accessibility"
$ws.Cells.Item(12,4).Value = "Describes the advantage of rescue robots being able to access tight or narrow spaces that are typically inaccessible to humans, enhancing their utility and safety in rescue operations."
$ws.Cells.Item(12,5).Value = 4
$ws.Cells.Item(12,6).Value = 30
# Row 13
$ws.Cells.Item(13,1).Value = "TL"
$ws.Cells.Item(13,2).Value = "perceived technological limitations"
$ws.Cells.Item(13,3).Value = "Perceived Risk"
$ws.Cells.Item(13,4).Value = "Captures instances where participants express skepticism about the ability of AI, such as rescue robots, to accurately assess complex situations compared to human judgment."
$ws.Cells.Item(13,6).Value = 13
# Row 14
$ws.Cells.Item(14,3).Value = "Adaptability"
$ws.Cells.Item(14,4).Value = "Captures instances where participants highlight the importance of rescue robots being able to adapt to various dangerous situations, such as extreme temperatures and physical conditions, in the context of discussing the perceived technological possibilities of rescue robots."
$ws.Cells.Item(14,6).Value = 4
# Row 15
$ws.Cells.Item(15,1).Value = "R"
$ws.Cells.Item(15,2).Value = "perceived risks"
$ws.Cells.Item(15,3).Value = "Technical Issues"
$ws.Cells.Item(15,4).Value = "Identifies instances where participants express concerns or frustrations related to the constraints or drawbacks of technical capabilities in rescue robots. This includes issues with limited senses, orientation, navigation, loss of control, and situations where autonomous systems fail to behave predictably due to software errors, as well as the inability to navigate complex environments or perform delicate tasks effectively."
$ws.Cells.Item(15,5).Value = 3
# Row 16
$ws.Cells.Item(16,1).Value = "SA"
$ws.Cells.Item(16,2).Value = "perceived safety"
$ws.Cells.Item(16,3).Value = "Focus on Task"
$ws.Cells.Item(16,4).Value = "Highlights rescue robots' efficiency and effectiveness in completing tasks without emotional interference, enhancing performance and reliability in rescue missions."
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 10
# Row 17
$ws.Cells.Item(17,1).Value = "TL"
$ws.Cells.Item(17,2).Value = "perceived technological limitations"
$ws.Cells.Item(17,3).Value = "Limited Flexibility"
$ws.Cells.Item(17,4).Value = "Highlights the robots' inability to adapt dynamically to unforeseen or complex situations."
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 2
# Row 18
$ws.Cells.Item(18,3).Value = "Specialized Tasks"
$ws.Cells.Item(18,4).Value = "Identifies instances where rescue robots are highlighted for their unique capabilities to perform specialized tasks beyond human capacity, such as flying, shrinking, hacking doors, and transmitting images for enhanced rescue operations."
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 10
# Row 19
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 8
# Row 20
$ws.Cells.Item(20,1).Value = "HRIP"
$ws.Cells.Item(20,2).Value = "perceived positive Human-Robot-Interaction"
$ws.Cells.Item(20,3).Value = "Collaborative Support"
$ws.Cells.Item(20,4).Value = "Identifies instances where rescue robots provide tangible support in rescue scenarios, including the delivery of essential resources and enhancing human capabilities through collaboration and support, rather than replacement."
$ws.Cells.Item(20,5).Value = 2
$ws.Cells.Item(20,6).Value = 2
# Row 21
$ws.Cells.Item(21,1).Value = "R"
$ws.Cells.Item(21,2).Value = "perceived risks"
$ws.Cells.Item(21,3).Value = "Potential Physical Harm"
$ws.Cells.Item(21,4).Value = "Identifies concerns related to the possibility of harm to individuals caused by technical malfunctions or unintended actions in rescue robots, such as malfunctions leading to accidental harm, injury, or even fatalities."
$ws.Cells.Item(21,5).Value = 2
$ws.Cells.Item(21,6).Value = 11
